# Add a new "october-2025" worksheet at the end of the workbook,
# mirroring the layout of the existing monthly sheets (single cell A1
# holding a "Total Gross Cumulative Voted Spending" summary string).

$wb = $excel.ActiveWorkbook

# Insert the new sheet immediately after the current last sheet so it
# lands at the end of the tab order (matching sheetId="22", rId="rId22").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "october-2025"

$newSheet.Range("A1").Value = " Total Gross Cumulative Voted Spending                               87,160       87,127          -33          0.0%               6,237         7.7%             107,015       103,472        3,543         3.4%"
